# Append one new row (row 65) to the portfolio data sheet, continuing the
# daily series for 2025-10-19 with the same SUZLON/TATAMOTORS/ETERNAL values
# as the prior day's row (row 64).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 65

# Force column A to be read as literal text (not auto-parsed into a date
# serial) while it is being populated, then drop the temporary "@" text
# format back to the sheet's default style so the new row doesn't pick up
# a stray per-cell style like the rest of the data rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-10-19"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 52.91999816894531
$ws.Cells.Item($row, 3).Value = 396.6000061035156
$ws.Cells.Item($row, 4).Value = 342.6499938964844
